$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14, Q1=15, matching style of O1 ---
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# --- Data rows 2-25: swap values in I/K/M/O, add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value2  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value2 = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value2 = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value2 = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value2 = 2   # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value2 = 2   # Q -> 2 (new)
}
